$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.444.14"
$ws.Range("E2").Value = "  -3.57%  "
$ws.Range("D3").Value = "2.476.64"
$ws.Range("E3").Value = "  -6.14%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'552.77"
$ws.Range("E5").Value = "  -4.74%  "
$ws.Range("D6").Value = "'146.58"
$ws.Range("E6").Value = "  -5.72%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.599"
$ws.Range("E8").Value = "  -3.69%  "
$ws.Range("D9").Value = "2.473.58"
$ws.Range("E9").Value = "  -6.14%  "
$ws.Range("E10").Value = "  -8.70%  "
$ws.Range("E11").Value = "  -5.89%  "
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "'0.356"
$ws.Range("E13").Value = "  -6.70%  "
$ws.Range("D14").Value = "'26.14"
$ws.Range("E14").Value = "  -7.85%  "
$ws.Range("D15").Value = "2.921.93"
$ws.Range("E15").Value = "  -6.19%  "
$ws.Range("E16").Value = "  -9.24%  "
$ws.Range("D17").Value = "61.321.87"
$ws.Range("E17").Value = "  -3.68%  "
$ws.Range("D18").Value = "2.475.83"
$ws.Range("E18").Value = "  -6.51%  "
$ws.Range("E19").Value = "  -7.79%  "
$ws.Range("D20").Value = "'7.05"
$ws.Range("E20").Value = "  -8.15%  "
$ws.Range("D21").Value = "'4.19"
$ws.Range("E21").Value = "  -7.43%  "
$ws.Range("D22").Value = "'321.60"
$ws.Range("E22").Value = "  -6.64%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'1.86"
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("D25").Value = "'64.02"
$ws.Range("E25").Value = "  -5.87%  "
$ws.Range("D26").Value = "0.0₃0992"
$ws.Range("E26").Value = "  -9.07%  "
$ws.Range("D27").Value = "2.606.13"
$ws.Range("E27").Value = "  -5.83%  "
$ws.Range("E28").Value = "  -6.06%  "
$ws.Range("D29").Value = "'542.47"
$ws.Range("E29").Value = "  -10.73%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").Value = "'8.33"
$ws.Range("E31").Value = "  -9.57%  "
$ws.Range("D32").Value = "'7.79"
$ws.Range("E32").Value = "  -3.88%  "
$ws.Range("E33").Value = "  -5.73%  "
$ws.Range("E34").Value = "  -8.20%  "
$ws.Range("E35").Value = "  -8.11%  "
$ws.Range("E36").Value = "  -10.01%  "
$ws.Range("D37").Value = "'4.88"
$ws.Range("E37").Value = "  -9.96%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'0.382"
$ws.Range("E39").Value = "  -5.08%  "
$ws.Range("D40").Value = "'18.56"
$ws.Range("E40").Value = "  -5.64%  "
$ws.Range("D41").Value = "'146.11"
$ws.Range("E41").Value = "  -3.21%  "
$ws.Range("D42").Value = "'1.73"
$ws.Range("E42").Value = "  -8.27%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "'40.43"
$ws.Range("E44").Value = "  -3.51%  "
$ws.Range("D45").Value = "'2.37"
$ws.Range("E45").Value = "  -7.25%  "
$ws.Range("D46").Value = "'147.53"
$ws.Range("E46").Value = "  -8.68%  "
$ws.Range("E47").Value = "  -6.44%  "
$ws.Range("D48").Value = "'21.20"
$ws.Range("E48").Value = "  -12.65%  "
$ws.Range("E49").Value = "  -7.63%  "
$ws.Range("D50").Value = "'0.597"
$ws.Range("E50").Value = "  -5.63%  "
$ws.Range("D51").Value = "'0.0945"
$ws.Range("E51").Value = "  -5.27%  "
